$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 and J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, centered, bordered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data values for columns I and J (rows 2-9)
$values = @{
    2 = 8
    3 = 4
    4 = 4
    5 = 7
    6 = 9
    7 = 7
    8 = 4
    9 = 4
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Cells.Item($row, 9).Value = $v
    $ws.Cells.Item($row, 10).Value = $v
}
